$d = $word.ActiveDocument

function Replace-Literal($findText, $replText, $wholeWord) {
    $r = $d.Content
    $found = $r.Find.Execute($findText, $true, $wholeWord, $false, $false, $false, $true, 1, $false)
    while ($found) {
        $r.Text = $replText
        $r.Collapse(0)
        $r.End = $d.Content.End
        $found = $r.Find.Execute($findText, $true, $wholeWord, $false, $false, $false, $true, 1, $false)
    }
}

Replace-Literal "Tatizo la mchwa - manukuu:" "The ants problem - subtitles:" $false
Replace-Literal "Mazungumzo huanza kwa sekunde 40 kwa hivyo niliongeza sekunde 27 kwa nyakati kama zilivyokuwa - John Argentino" "The dialogue starts at 40 seconds in so I added 27 seconds to the times as they were - John Argentino" $false
Replace-Literal "[Muziki]" "[Music]" $false
Replace-Literal "sawa kwa hivyo mafumbo nitaenda" "okay so the puzzles I'm going to" $false
Replace-Literal "changamoto uliyonayo ni mbili za msingi" "challenge you with are two basic" $false
Replace-Literal "matoleo ya fumbo ngumu zaidi" "versions of a more complicated puzzle" $false
Replace-Literal "inayojulikana kama fumbo la mchwa, ambalo mimi ni" "known as the ants puzzle, which I'm" $false
Replace-Literal "pengine kwenda kujadili katika tofauti" "probably going to discuss in a different" $false
Replace-Literal "video. Ngoja nimalizie kuandika" "video. Let me just finish writing down" $false
Replace-Literal "kichwa na, vizuri, naweza hata kuchora a" "the title and, well, I can even draw a" $false
Replace-Literal "mchwa mdogo hapa. sawa, tupate" "little ant right here. okay, let's get" $false
Replace-Literal "imeanza! Kama nilivyosema nitajadili" "started! As I said I'm going to discuss" $false
Replace-Literal "mafumbo mawili katika fumbo la kwanza hapo" "two puzzles in the first puzzle there" $false
Replace-Literal "ni mchwa wawili kwenye kinyesi cha juu sana: aina" "are two ants on a very high stool: a sort" $false
Replace-Literal "ya Mlima, gorofa juu na mbili" "of Mountain, flat on the top with two" $false
Replace-Literal "miamba mikali kwa pande zote mbili. Gorofa" "steep cliffs to both the sides. The flat" $false
Replace-Literal "kilele ni mita moja upana wa mchwa wawili hoja" "peak is one meter wide the two ants move" $false
Replace-Literal "kwa kasi, tuiite V, ambayo ni" "with a velocity, let's call it V, which is" $false
Replace-Literal "sawa kwa wote wawili na hiyo ni" "the same for both of them and that is" $false
Replace-Literal "sawa na sentimita moja kwa sekunde. Wewe" "equal to one centimeter per second. You" $false
Replace-Literal "inaweza kuamua mwelekeo kuelekea kila mmoja" "can decide the direction towards each" $false
Replace-Literal "mchwa husogea ikiwa ni kulia au kushoto na" "ant moves if it is right or left and" $false
Replace-Literal "wapi hasa kuweka mchwa wawili kwenye" "where exactly to place the two ants on the" $false
Replace-Literal "juu ya mlima. Kusudi lako ni" "top of the mountain. Your purpose is to" $false
Replace-Literal "fanya wakati mchwa wa mwisho huchukua hapo awali" "make the time the last ant takes before" $false
Replace-Literal "kuanguka kwa muda mrefu iwezekanavyo. Mchwa hawawezi" "falling the longest possible. Ants cannot" $false
Replace-Literal "tulia: lazima wahamie kulia au" "be still: they must move to the right or" $false
Replace-Literal "upande wa kushoto lakini lazima wasogee na baada" "to the left but they must move and after" $false
Replace-Literal "wakikutana wanageuka na" "meeting each other they turn around and" $false
Replace-Literal "endelea kusonga na sawa lakini kinyume" "keep moving with the same but opposite" $false
Replace-Literal "kwa hivyo tena ni nafasi gani sahihi" "so again what are the precise positions" $false
Replace-Literal "ambapo ninapaswa kuwaweka mchwa wawili ndani" "where I should place the two ants in" $false
Replace-Literal "ili kupata muda mrefu zaidi kabla ya" "order to get the longest time before the" $false
Replace-Literal "chungu mwisho huanguka? Fumbo la pili ni" "last ant falls? The second puzzle is" $false
Replace-Literal "kimsingi ni sawa lakini sasa tuna tatu" "basically the same but now we have three" $false
Replace-Literal "mchwa badala ya wawili." "ants instead of two." $false
Replace-Literal "Kama kabla ya mchwa kasi ni moja" "As before the ants velocity is one" $false
Replace-Literal "sentimita kwa sekunde, kila mchwa hugeuka" "centimeter per second, every ant turns" $false
Replace-Literal "karibu baada ya kukutana na mchwa mwingine na" "around after meeting another ant and" $false
Replace-Literal "kilele kina upana wa mita moja. Hivyo, ni nini" "the peak is one meter wide. So, what are" $false
Replace-Literal "sasa nafasi sahihi" "now the precise positions" $false
Replace-Literal "Ninapaswa kuweka mchwa watatu kwa mpangilio" "I should place the three ants in order" $false
Replace-Literal "kupata muda mrefu zaidi kabla ya mwisho" "to get the longest time before the last" $false
Replace-Literal "chungu huanguka chini? Natumaini ulifurahia hili" "ant falls down? I hope you enjoyed this" $false
Replace-Literal "video fanya bora na bahati nzuri" "video do your best and good luck" $false
Replace-Literal "kasi" "velocity" $true
